# The "yesNo(boolean) ... failed" sample stack trace embedded in the
# document reflects the Acceleo / M2Doc / JUnit source line numbers of
# the versions in use when the expected-generation fixture was last
# (re)generated. Moving from 3.2.0 to 3.2.1 shifted a number of those
# line numbers and also changed the shape of the nested JUnit
# ParentRunner frames near the bottom of the trace. Replace the whole
# bold/red stack-trace run in one shot with the updated text.
$d = $word.ActiveDocument

$oldText = @(
    "yesNo(boolean) with arguments [null] failed:",
    "`tnull",
    "java.lang.IllegalArgumentException",
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)",
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)",
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "`tat java.lang.reflect.Method.invoke(Method.java:498)",
    "`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)",
    "`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)",
    "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)",
    "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)",
    "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208)",
    "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)",
    "`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)",
    "`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:591)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)",
    "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1556)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)",
    "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:301)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)",
    "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:286)",
    "`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)",
    "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)",
    "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)",
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)",
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)",
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "`tat java.lang.reflect.Method.invoke(Method.java:498)",
    "`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)",
    "`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)",
    "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)",
    "`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)",
    "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)",
    "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)",
    "`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)",
    "`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)",
    "`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)",
    "`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)",
    "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)",
    "`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:128)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:27)",
    "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)",
    "`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)",
    "`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)",
    "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)",
    "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:128)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:27)",
    "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)",
    "`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)",
    "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:128)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:27)",
    "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)",
    "`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)",
    "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)",
    "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)",
    "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
) -join "`n"
$oldText = $oldText + "`n"

$newText = @(
    "yesNo(boolean) with arguments [null] failed:",
    "`tnull",
    "java.lang.IllegalArgumentException",
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)",
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)",
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "`tat java.lang.reflect.Method.invoke(Method.java:498)",
    "`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)",
    "`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)",
    "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)",
    "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)",
    "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208)",
    "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:192)",
    "`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)",
    "`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:604)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)",
    "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1675)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)",
    "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)",
    "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299)",
    "`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)",
    "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)",
    "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)",
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)",
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)",
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)",
    "`tat java.lang.reflect.Method.invoke(Method.java:498)",
    "`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)",
    "`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)",
    "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)",
    "`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)",
    "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)",
    "`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)",
    "`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)",
    "`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)",
    "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)",
    "`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:128)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:27)",
    "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)",
    "`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)",
    "`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)",
    "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:128)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:27)",
    "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)",
    "`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:128)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:27)",
    "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)",
    "`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)",
    "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)",
    "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
) -join "`n"
$newText = $newText + "`n"

$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) {
    throw "Could not find the stack trace text to replace"
}
